# Updated cryptos list on Sat Jul 20 11:41:21 UTC 2024 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns for each coin row, and
# swaps the OKB / InjectiveProtocol / Maker rows (44-46) into their new order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For price cells whose new value looks like a plain decimal number (e.g.
# "590.96"), force the cell to Text format first so Excel keeps it as a
# literal string "590.96" instead of silently converting it to the numeric
# value 590.96 (the source data intentionally stores these prices as text,
# as proven by the many "thousands.hundreds" values such as "66.548.35"
# that can never be parsed as a single number).
$ws.Range('D2').Value = '66.548.35'
$ws.Range('E2').Value = '  +3.99%  '
$ws.Range('D3').Value = '3.490.68'
$ws.Range('E3').Value = '  +2.67%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '590.96'
$ws.Range('E5').Value = '  +3.27%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '168.51'
$ws.Range('E6').Value = '  +3.90%  '
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').Value = '3.487.76'
$ws.Range('E8').Value = '  +2.59%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.591'
$ws.Range('E9').Value = '  +7.53%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.33'
$ws.Range('E10').Value = '  +0.44%  '
$ws.Range('E11').Value = '  +6.21%  '
$ws.Range('E12').Value = '  +3.43%  '
$ws.Range('D13').Value = '4.092.77'
$ws.Range('E13').Value = '  +2.53%  '
$ws.Range('E14').Value = '  -0.50%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '28.08'
$ws.Range('E15').Value = '  +4.87%  '
$ws.Range('D16').Value = '66.571.11'
$ws.Range('E16').Value = '  +3.97%  '
$ws.Range('E17').Value = '  +2.92%  '
$ws.Range('D18').Value = '3.483.89'
$ws.Range('E18').Value = '  +1.81%  '
$ws.Range('E19').Value = '  +3.04%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.94'
$ws.Range('E20').Value = '  +4.21%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '390.40'
$ws.Range('E21').Value = '  +4.69%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.91'
$ws.Range('E22').Value = '  +1.91%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '72.82'
$ws.Range('E23').Value = '  +3.62%  '
$ws.Range('E24').Value = '  +0.41%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.532'
$ws.Range('E25').Value = '  +4.09%  '
$ws.Range('E26').Value = '  +6.49%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.22'
$ws.Range('E27').Value = '  +8.37%  '
$ws.Range('E28').Value = '  +1.52%  '
$ws.Range('E29').Value = '  +0.35%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.31'
$ws.Range('E30').Value = '  +4.61%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.45'
$ws.Range('E31').Value = '  +4.83%  '
$ws.Range('E32').Value = '  +2.87%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '23.54'
$ws.Range('E33').Value = '  +3.63%  '
$ws.Range('E34').Value = '  +4.75%  '
$ws.Range('E35').Value = '  +0.06%  '
$ws.Range('E36').Value = '  +8.99%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '162.62'
$ws.Range('E37').Value = '  +2.44%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.888'
$ws.Range('E38').Value = '  +4.09%  '
$ws.Range('E39').Value = '  +5.64%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.76'
$ws.Range('E40').Value = '  +5.84%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0743'
$ws.Range('E41').Value = '  +2.90%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.63'
$ws.Range('E42').Value = '  +6.88%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '26.35'
$ws.Range('E43').Value = '  +2.46%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '2.781.54'
$ws.Range('E44').Value = '  +2.16%  '
$ws.Range('B45').Value = 'OKB'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '43.06'
$ws.Range('E45').Value = '  +1.04%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '26.57'
$ws.Range('E46').Value = '  +3.19%  '
$ws.Range('E47').Value = '  +2.53%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.49'
$ws.Range('E48').Value = '  +4.28%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '345.27'
$ws.Range('E49').Value = '  +5.52%  '
$ws.Range('E50').Value = '  +4.57%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '33.72'
$ws.Range('E51').Value = '  +12.22%  '
